$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old header row (row 1) with the new, renamed headers.
$ws.Range("A1").Value = "DefFirstName"
$ws.Range("B1").Value = "CaseNumber"
$ws.Range("C1").Value = "SubCaseNumber"
$ws.Range("D1").Value = "DefLastName"
$ws.Range("E1").Value = "ChargeDescription"
$ws.Range("F1").Value = "SectionCode"
$ws.Range("G1").Value = "DegreeCode"
$ws.Range("H1").Value = "InsuranceStatus"
$ws.Range("I1").Value = "IsMoving"
$ws.Range("J1").Value = "AttorneyLastName"
$ws.Range("K1").Value = "AttorneyFirstName"
$ws.Range("L1").Value = "PubDef"

# New insurance-status value for row 2.
$ws.Range("H2").Value = "U"

# Column I ("IsMoving") was stored as boolean; rewrite as plain numbers
# (0/1) so the underlying cell type switches from boolean to numeric,
# keeping the same logical value.
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("I18").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("I22").Value = 1

# Column L ("PubDef") gains new numeric values for every data row.
$ws.Range("L2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("L16").Value = 1
$ws.Range("L17").Value = 1
$ws.Range("L18").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("L21").Value = 1
$ws.Range("L22").Value = 0

# Restore the scroll position / selection that was active when the
# workbook was last saved.
$ws.Range("I23").Select()
